$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.140.17'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '3.534.32'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.46%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.140'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '4.145.29'
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +10.32%  '
$ws.Range("E14").Value = '  +1.05%  '
$ws.Range("D15").Value = '68.082.98'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000181'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '3.537.84'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '401.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '24.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  -2.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.883'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("E40").Value = '  +6.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.933.90'
$ws.Range("E42").Value = '  +3.23%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0741'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("E46").Value = '  -1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '352.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0307'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("E51").Value = '  +2.98%  '
